# Updated cryptos list on Wed Feb 14 14:26:30 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'51.594.03"
$ws.Range("E2").Value = "  +5.91%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'2.757.50"
$ws.Range("E3").Value = "  +5.82%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.41%  "

# Row 5 - Solana
$ws.Range("D5").Value = "'117.56"
$ws.Range("E5").Value = "  +7.21%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'332.63"
$ws.Range("E6").Value = "  +4.06%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.534"
$ws.Range("E7").Value = "  +3.35%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.16%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +7.22%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "'41.45"
$ws.Range("E10").Value = "  +5.87%  "

# Row 11 - Chainlink
$ws.Range("D11").Value = "'20.23"
$ws.Range("E11").Value = "  +3.03%  "

# Row 12 - Dogecoin
$ws.Range("D12").Value = "'0.0831"
$ws.Range("E12").Value = "  +3.69%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +3.35%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +6.74%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "'3.200.37"
$ws.Range("E15").Value = "  +5.03%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "'2.765.98"
$ws.Range("E16").Value = "  +4.71%  "

# Row 17 - Polygon
$ws.Range("D17").Value = "'0.885"
$ws.Range("E17").Value = "  +4.31%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "'51.604.52"
$ws.Range("E18").Value = "  +5.46%  "

# Row 19 - InternetComputer(DFINITY)
$ws.Range("D19").Value = "'13.70"
$ws.Range("E19").Value = "  +7.59%  "

# Row 20 - ImmutableX
$ws.Range("E20").Value = "  +4.76%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +3.85%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "'0.0₃0964"
$ws.Range("E22").Value = "  +3.31%  "

# Row 23 - BitcoinCash
$ws.Range("D23").Value = "'277.07"
$ws.Range("E23").Value = "  +3.10%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'70.16"
$ws.Range("E24").Value = "  +0.70%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +5.97%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "'26.87"
$ws.Range("E26").Value = "  +3.18%  "

# Row 27 - LEO
$ws.Range("D27").Value = "'4.14"
$ws.Range("E27").Value = "  +1.03%  "

# Row 28 - Dai
$ws.Range("E28").Value = "  -0.06%  "

# Row 29 - Cosmos
$ws.Range("D29").Value = "'10.29"
$ws.Range("E29").Value = "  +2.87%  "

# Row 30 - Toncoin
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = "  +0.27%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +3.47%  "

# Row 32 - InjectiveProtocol
$ws.Range("D32").Value = "'35.61"
$ws.Range("E32").Value = "  +1.47%  "

# Row 33 - OKB
$ws.Range("D33").Value = "'50.49"
$ws.Range("E33").Value = "  +2.47%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +4.16%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "'0.0831"
$ws.Range("E35").Value = "  +5.68%  "

# Row 36 - Celestia
$ws.Range("D36").Value = "'19.38"
$ws.Range("E36").Value = "  +2.45%  "

# Row 37/38 swapped: ARBITRUM <-> FirstDigitalUSD
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'2.11"
$ws.Range("E38").Value = "  +5.52%  "

# Row 39 - RenderToken
$ws.Range("E39").Value = "  +2.61%  "

# Row 40 - LidoDAOToken
$ws.Range("D40").Value = "'3.26"
$ws.Range("E40").Value = "  +6.26%  "

# Row 41 - Monero
$ws.Range("D41").Value = "'130.66"
$ws.Range("E41").Value = "  +5.24%  "

# Row 42 - EnergySwap
$ws.Range("D42").Value = "'23.41"
$ws.Range("E42").Value = "  +6.32%  "

# Row 43/44 swapped: VeChain <-> Stellar
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").Value = "'0.114"
$ws.Range("E43").Value = "  +3.54%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0344"
$ws.Range("E44").Value = "  +10.32%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  +4.06%  "

# Row 46 - Stacks
$ws.Range("E46").Value = "  +13.42%  "

# Row 47 - Maker
$ws.Range("D47").Value = "'2.117.30"
$ws.Range("E47").Value = "  +2.26%  "

# Row 48 - NEARProtocol
$ws.Range("D48").Value = "'3.35"
$ws.Range("E48").Value = "  +5.16%  "

# Row 49 - ApeXProtocol
$ws.Range("E49").Value = "  +2.78%  "

# Row 50 - THORChain
$ws.Range("D50").Value = "'5.61"
$ws.Range("E50").Value = "  +8.54%  "

# Row 51 - FraxShare
$ws.Range("D51").Value = "'9.02"
$ws.Range("E51").Value = "  +2.47%  "
